# Fix Training Data Issue
# The BF column holds a "Date" label (row 1) followed by a date value per
# team row (rows 2-31). The stored value "5-3-2013-14" was off by one day
# due to how NBA stats were displayed; replace it with the correct date
# "2014-05-03" for every data row, keeping the cell as plain text (not an
# auto-converted date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-3-2013-14"
$newValue = "2014-05-03"

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 31) { $lastRow = 31 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 58)  # column BF
    if ($cell.Value2 -eq $oldValue) {
        # Force the value to be stored as text rather than letting Excel's
        # smart-parser reinterpret the ISO-like string as a date serial.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.ClearFormats()
    }
}
